$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(243, 243, 8, 1, 38, "09:27:35"),
    @(282, 243, 13, 21, 28, "14:47:35"),
    @(360, 243, 7, 57, 26, "09:23:50"),
    @(399, 243, 9, 5, 33, "10:32:06"),
    @(263, 263, 9, 3, 56, "10:42:55"),
    @(380, 263, 0, 24, 54, "02:07:37"),
    @(419, 263, 0, 24, 24, "02:07:16"),
    @(458, 263, 8, 50, 38, "10:33:39"),
    @(497, 263, 8, 50, 14, "10:33:25"),
    @(536, 263, 2, 43, 58, "04:27:18"),
    @(575, 263, 1, 58, 25, "03:41:54"),
    @(614, 263, 1, 28, 36, "03:12:14"),
    @(653, 263, 13, 34, 57, "15:18:44"),
    @(731, 263, 12, 11, 48, "13:55:50"),
    @(809, 263, 0, 50, 39, "02:34:55"),
    @(926, 263, 15, 21, 38, "17:06:16"),
    @(965, 263, 15, 52, 12, "17:36:59"),
    @(1004, 263, 9, 39, 14, "11:24:09")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $arr = New-Object "object[,]" 1,6
    for ($j = 0; $j -lt 6; $j++) {
        $arr[0,$j] = $rowVals[$j]
    }
    $r = $i + 1
    $rangeAddr = "A" + $r + ":F" + $r
    $ws.Range($rangeAddr).Value = $arr
}

$ws.Range("A19:F20").Clear()